$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "2025/12/03 07:00"
$ws.Range("B16").Value = "43,758位本"
$ws.Range("C16").Value = "97位 広告・宣伝 (本)"
$ws.Range("D16").Value = "175位商業デザイン"
$ws.Range("E16").Value = "2,145位ビジネス実用本"
$ws.Range("F16").Value = "-"
$ws.Range("G16").Value = "-"
